$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Cilantro needs to be inserted at row 353,
# pushing the existing rows 353-378 down to 354-379 (dimension A1:R378 -> A1:R379).
$ws.Rows(353).Insert()

$ws.Cells.Item(353, 1).Value = 3
$ws.Cells.Item(353, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(353, 3).Value = "Coquimbo"
$ws.Cells.Item(353, 4).Value = 44746
$ws.Cells.Item(353, 5).Value = 5
$ws.Cells.Item(353, 6).Value = 100112040
$ws.Cells.Item(353, 7).Value = "Cilantro"
$ws.Cells.Item(353, 8).Value = "Sin especificar"
$ws.Cells.Item(353, 9).Value = "Primera"
$ws.Cells.Item(353, 10).Value = 172
$ws.Cells.Item(353, 11).Value = 3500
$ws.Cells.Item(353, 12).Value = 4000
$ws.Cells.Item(353, 13).Value = 3753
$ws.Cells.Item(353, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(353, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(353, 16).Value = 1251
$ws.Cells.Item(353, 17).Value = 3
$ws.Cells.Item(353, 18).Value = "Hortaliza"
